$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(99, 8).Value = 1376.625  # H99
$ws.Cells.Item(99, 9).Value = 530.2  # I99
$ws.Cells.Item(99, 10).Value = 2787.3333  # J99
$ws.Cells.Item(99, 11).Value = 1590.6  # K99
$ws.Cells.Item(99, 12).Value = 8361.999899999999  # L99
$ws.Cells.Item(99, 13).Value = -92.60000000000014  # M99
$ws.Cells.Item(99, 14).Value = -11357.9999  # N99
$ws.Cells.Item(100, 8).Value = 720.8  # H100
$ws.Cells.Item(100, 10).Value = 553  # J100
$ws.Cells.Item(100, 12).Value = 553  # L100
$ws.Cells.Item(100, 14).Value = -1635  # N100
$ws.Cells.Item(123, 8).Value = 37156.25  # H123
$ws.Cells.Item(123, 10).Value = 37156.25  # J123
$ws.Cells.Item(123, 12).Value = 37156.25  # L123
$ws.Cells.Item(123, 14).Value = -46956.25  # N123
$ws.Cells.Item(125, 8).Value = 1735.1538  # H125
$ws.Cells.Item(125, 9).Value = 499.125  # I125
$ws.Cells.Item(125, 10).Value = 3712.8  # J125
$ws.Cells.Item(125, 11).Value = 4492.125  # K125
$ws.Cells.Item(125, 12).Value = 33415.2  # L125
$ws.Cells.Item(125, 13).Value = -2032.125  # M125
$ws.Cells.Item(125, 14).Value = -38335.2  # N125
$ws.Cells.Item(129, 8).Value = 833.2963  # H129
$ws.Cells.Item(129, 10).Value = 955.9761999999999  # J129
$ws.Cells.Item(129, 12).Value = 2867.9286  # L129
$ws.Cells.Item(129, 14).Value = -12867.9286  # N129
$ws.Cells.Item(132, 8).Value = 1258747.5  # H132
$ws.Cells.Item(132, 9).Value = 2342  # I132
$ws.Cells.Item(132, 10).Value = 24502250  # J132
$ws.Cells.Item(132, 11).Value = 7026  # K132
$ws.Cells.Item(132, 12).Value = 73506750  # L132
$ws.Cells.Item(132, 13).Value = -4496  # M132
$ws.Cells.Item(132, 14).Value = -73511810  # N132
$ws.Cells.Item(138, 8).Value = 2318779.5  # H138
$ws.Cells.Item(138, 9).Value = 1949.2413  # I138
$ws.Cells.Item(138, 10).Value = 3881293  # J138
$ws.Cells.Item(138, 11).Value = 5847.7239  # K138
$ws.Cells.Item(138, 12).Value = 11643879  # L138
$ws.Cells.Item(138, 13).Value = -707.7239  # M138
$ws.Cells.Item(138, 14).Value = -11654159  # N138

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4347283  # H32
$ws.Cells.Item(32, 9).Value = 5309419  # I32
$ws.Cells.Item(32, 11).Value = 5309419  # K32
$ws.Cells.Item(32, 13).Value = -5309132  # M32
$ws.Cells.Item(74, 8).Value = 5861405.5  # H74
$ws.Cells.Item(74, 9).Value = 7606767.5  # I74
$ws.Cells.Item(74, 10).Value = 101711  # J74
$ws.Cells.Item(74, 11).Value = 7606767.5  # K74
$ws.Cells.Item(74, 12).Value = 101711  # L74
$ws.Cells.Item(74, 13).Value = -7605893.5  # M74
$ws.Cells.Item(74, 14).Value = -103459  # N74
$ws.Cells.Item(77, 8).Value = 5861405.5  # H77
$ws.Cells.Item(77, 9).Value = 7606767.5  # I77
$ws.Cells.Item(77, 10).Value = 101711  # J77
$ws.Cells.Item(77, 11).Value = 38033837.5  # K77
$ws.Cells.Item(77, 12).Value = 508555  # L77
$ws.Cells.Item(77, 13).Value = -38029469.5  # M77
$ws.Cells.Item(77, 14).Value = -517291  # N77
$ws.Cells.Item(122, 8).Value = 37042044  # H122
$ws.Cells.Item(122, 9).Value = 7512  # I122
$ws.Cells.Item(122, 11).Value = 22536  # K122
$ws.Cells.Item(122, 13).Value = -20086  # M122
$ws.Cells.Item(132, 8).Value = 225555.33  # H132
$ws.Cells.Item(132, 9).Value = 253000  # I132
$ws.Cells.Item(132, 10).Value = 203599.6  # J132
$ws.Cells.Item(132, 11).Value = 759000  # K132
$ws.Cells.Item(132, 12).Value = 610798.8  # L132
$ws.Cells.Item(132, 13).Value = -756470  # M132
$ws.Cells.Item(132, 14).Value = -615858.8  # N132

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1309  # H94
$ws.Cells.Item(94, 9).Value = 500  # I94
$ws.Cells.Item(94, 11).Value = 500  # K94
$ws.Cells.Item(94, 13).Value = -49  # M94

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 25642640  # H16
$ws.Cells.Item(16, 9).Value = 1778.091  # I16
$ws.Cells.Item(16, 11).Value = 1778.091  # K16
$ws.Cells.Item(16, 13).Value = -1491.091  # M16
$ws.Cells.Item(106, 8).Value = 11144.615  # H106
$ws.Cells.Item(106, 10).Value = 11144.615  # J106
$ws.Cells.Item(106, 12).Value = 11144.615  # L106
$ws.Cells.Item(106, 14).Value = -13668.615  # N106
$ws.Cells.Item(113, 8).Value = 25642640  # H113
$ws.Cells.Item(113, 9).Value = 1778.091  # I113
$ws.Cells.Item(113, 11).Value = 1778.091  # K113
$ws.Cells.Item(113, 13).Value = 391.9090000000001  # M113

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 35082.137  # H5
$ws.Cells.Item(5, 9).Value = 71914.86  # I5
$ws.Cells.Item(5, 10).Value = 704.93335  # J5
$ws.Cells.Item(5, 11).Value = 215744.58  # K5
$ws.Cells.Item(5, 12).Value = 2114.80005  # L5
$ws.Cells.Item(5, 13).Value = -215632.58  # M5
$ws.Cells.Item(5, 14).Value = -2338.80005  # N5
$ws.Cells.Item(107, 8).Value = 869.5223999999999  # H107
$ws.Cells.Item(107, 9).Value = 413.29788  # I107
$ws.Cells.Item(107, 10).Value = 1941.65  # J107
$ws.Cells.Item(107, 11).Value = 1239.89364  # K107
$ws.Cells.Item(107, 12).Value = 5824.950000000001  # L107
$ws.Cells.Item(107, 13).Value = 680.10636  # M107
$ws.Cells.Item(107, 14).Value = -9664.950000000001  # N107
$ws.Cells.Item(113, 8).Value = 516.9828  # H113
$ws.Cells.Item(113, 9).Value = 468.4  # I113
$ws.Cells.Item(113, 10).Value = 542.5526  # J113
$ws.Cells.Item(113, 11).Value = 1405.2  # K113
$ws.Cells.Item(113, 12).Value = 1627.6578  # L113
$ws.Cells.Item(113, 13).Value = 764.8000000000002  # M113
$ws.Cells.Item(113, 14).Value = -5967.6578  # N113
$ws.Cells.Item(131, 8).Value = 803.7436  # H131
$ws.Cells.Item(131, 9).Value = 473  # I131
$ws.Cells.Item(131, 10).Value = 917.7931  # J131
$ws.Cells.Item(131, 11).Value = 1419  # K131
$ws.Cells.Item(131, 12).Value = 2753.3793  # L131
$ws.Cells.Item(131, 13).Value = 3621  # M131
$ws.Cells.Item(131, 14).Value = -12833.3793  # N131
$ws.Cells.Item(135, 8).Value = 35082.137  # H135
$ws.Cells.Item(135, 9).Value = 71914.86  # I135
$ws.Cells.Item(135, 10).Value = 704.93335  # J135
$ws.Cells.Item(135, 11).Value = 647233.74  # K135
$ws.Cells.Item(135, 12).Value = 6344.40015  # L135
$ws.Cells.Item(135, 13).Value = -644698.74  # M135
$ws.Cells.Item(135, 14).Value = -11414.40015  # N135

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1816.7858  # H113
$ws.Cells.Item(113, 9).Value = 1635.1666  # I113
$ws.Cells.Item(113, 10).Value = 1953  # J113
$ws.Cells.Item(113, 11).Value = 1635.1666  # K113
$ws.Cells.Item(113, 12).Value = 1953  # L113
$ws.Cells.Item(113, 13).Value = 534.8334  # M113
$ws.Cells.Item(113, 14).Value = -6293  # N113

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 2259.2307  # H93
$ws.Cells.Item(93, 9).Value = 2234.7273  # I93
$ws.Cells.Item(93, 11).Value = 2234.7273  # K93
$ws.Cells.Item(93, 13).Value = -986.7273  # M93
$ws.Cells.Item(122, 8).Value = 3499.6667  # H122
$ws.Cells.Item(122, 9).Value = 3747  # I122
$ws.Cells.Item(122, 10).Value = 3005  # J122
$ws.Cells.Item(122, 11).Value = 11241  # K122
$ws.Cells.Item(122, 12).Value = 9015  # L122
$ws.Cells.Item(122, 13).Value = -8791  # M122
$ws.Cells.Item(122, 14).Value = -13915  # N122
$ws.Cells.Item(132, 8).Value = 49025.816  # H132
$ws.Cells.Item(132, 9).Value = 2660.3333  # I132
$ws.Cells.Item(132, 11).Value = 7980.999899999999  # K132
$ws.Cells.Item(132, 13).Value = -5450.999899999999  # M132

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 0  # H46
$ws.Cells.Item(46, 10).Value = 0  # J46
$ws.Cells.Item(46, 12).Value = 0  # L46
$ws.Cells.Item(46, 14).ClearContents()  # N46
$ws.Cells.Item(124, 8).Value = 0  # H124
$ws.Cells.Item(124, 10).Value = 0  # J124
$ws.Cells.Item(124, 12).Value = 0  # L124
$ws.Cells.Item(124, 14).ClearContents()  # N124
$ws.Cells.Item(132, 8).Value = 47131.953  # H132
$ws.Cells.Item(132, 9).Value = 38598.52  # I132
$ws.Cells.Item(132, 10).Value = 60685.06  # J132
$ws.Cells.Item(132, 11).Value = 115795.56  # K132
$ws.Cells.Item(132, 12).Value = 182055.18  # L132
$ws.Cells.Item(132, 13).Value = -113265.56  # M132
$ws.Cells.Item(132, 14).Value = -187115.18  # N132
$ws.Cells.Item(134, 8).Value = 0  # H134
$ws.Cells.Item(134, 10).Value = 0  # J134
$ws.Cells.Item(134, 12).Value = 0  # L134
$ws.Cells.Item(134, 14).ClearContents()  # N134
